$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update quantities: Ordinateur portable (B3) and Imprimante HP (B6)
$ws.Range("B3").Value = 16
$ws.Range("B6").Value = 116
